$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.246.31'
$ws.Range("E2").Value = '  +2.68%  '

# Row 3
$ws.Range("D3").Value = '1.718.09'
$ws.Range("E3").Value = '  +3.00%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4708'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.95%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2622'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.64%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06191'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.19%  '

# Row 10
$ws.Range("D10").Value = '1.715.50'
$ws.Range("E10").Value = '  +3.30%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07064'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.68%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5956'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.81%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.421'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.98%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.21%  '

# Row 16
$ws.Range("E16").Value = '  +0.06%  '

# Row 17
$ws.Range("E17").Value = '  +0.09%  '

# Row 18
$ws.Range("D18").Value = '26.257.73'
$ws.Range("E18").Value = '  +2.77%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006796'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.45%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.24%  '

# Row 21
$ws.Range("D21").Value = '1.936.07'
$ws.Range("E21").Value = '  +3.32%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.531'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.701'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.30%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.240'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.96%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.85%  '

# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.761'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.73%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.398'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.01%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.59%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.943'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.71%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.681'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.78%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07749'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.21%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04493'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.92%  '

# Row 34
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.000'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.12%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.616'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.72%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9741'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.25%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6178'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.64%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9307'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.43%  '

# Row 39
$ws.Range("B39").Value = 'Quant'
$ws.Range("C39").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '114.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +17.60%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.438'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.12%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.921'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.62%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.622'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +15.93%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01479'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.74%  '

# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3820'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.27%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1176'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.70%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.292'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.13%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05263'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.777'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.97%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.89%  '

# Row 51
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3375'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.19%  '
